# Fixed a bug in Respin
# Rewrites the Respin data table (A2:F20) with the corrected row order/values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(901,  16, 15, 45, 60, 60),
    @(1001, 18, 30, 75, 60, 72),
    @(301,   6, 45, 30, 60, 45),
    @(501,   9, 52, 30, 75, 45),
    @(201,   9, 30, 15, 45, 30),
    @(1201,  2, 10, 10, 10, 10),
    @(902,   1,  0,  0,  0,  0),
    @(701,   3, 90, 45, 97, 15),
    @(601,   9, 60, 67, 60, 42),
    @(101,   9, 30, 15, 60, 15),
    @(401,   9, 48, 67, 75, 45),
    @(801,   3, 67, 65, 52, 45),
    @(1202,  2, 10, 10, 10, 10),
    @(1203,  3, 15, 15, 15, 15),
    @(1,     0,  2,  2,  2,  2),
    @(1101,  0, 15, 30, 30,  0),
    @(2,     0,  2,  2,  2,  2),
    @(502,   0,  4,  0,  0,  0),
    @(802,   0,  4,  5,  4,  0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($startRow + $i, $j + 1).Value = $row[$j]
    }
}
